$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Insert the new "metadata" sheet right after "data" so tab order is data, metadata
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row (row 1) - columns B..G
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Miscellaneous Metabolic Disorders"
$ws.Range("C2").Value = 3468
$ws.Range("D2").Value = "1.8"
$ws.Range("E2").Value = "2021-08-29T07:03:32.153926Z"
$ws.Range("F2").Value = "2021-10-05 14:34:53.040487"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3468/?format=json"

# Match the header style used on the "data" sheet (bold, thin border, centered)
$headerStyleRange = $dataSheet.Range("B1")
$ws.Range("B1:G1").Font.Bold = $headerStyleRange.Font.Bold
$ws.Range("B1:G1").HorizontalAlignment = $headerStyleRange.HorizontalAlignment
$ws.Range("B1:G1").VerticalAlignment = $headerStyleRange.VerticalAlignment
$ws.Range("B1:G1").Borders.LineStyle = $headerStyleRange.Borders.LineStyle
$ws.Range("A2").Font.Bold = $dataSheet.Range("A2").Font.Bold
$ws.Range("A2").HorizontalAlignment = $dataSheet.Range("A2").HorizontalAlignment
$ws.Range("A2").VerticalAlignment = $dataSheet.Range("A2").VerticalAlignment
$ws.Range("A2").Borders.LineStyle = $dataSheet.Range("A2").Borders.LineStyle

# Update the "time_taken" column (F) on the "data" sheet with refreshed timestamps
$timestamps = @(
  "2021-10-05 14:34:53.043669","2021-10-05 14:34:53.043676","2021-10-05 14:34:53.043679","2021-10-05 14:34:53.043682",
  "2021-10-05 14:34:53.043685","2021-10-05 14:34:53.043687","2021-10-05 14:34:53.043690","2021-10-05 14:34:53.043692",
  "2021-10-05 14:34:53.043695","2021-10-05 14:34:53.043698","2021-10-05 14:34:53.043700","2021-10-05 14:34:53.043702",
  "2021-10-05 14:34:53.043705","2021-10-05 14:34:53.043707","2021-10-05 14:34:53.043709","2021-10-05 14:34:53.043712",
  "2021-10-05 14:34:53.043714","2021-10-05 14:34:53.043717","2021-10-05 14:34:53.043719","2021-10-05 14:34:53.043722",
  "2021-10-05 14:34:53.043724","2021-10-05 14:34:53.043727","2021-10-05 14:34:53.043729","2021-10-05 14:34:53.043731",
  "2021-10-05 14:34:53.043734","2021-10-05 14:34:53.043737","2021-10-05 14:34:53.043739","2021-10-05 14:34:53.043742",
  "2021-10-05 14:34:53.043744","2021-10-05 14:34:53.043746","2021-10-05 14:34:53.043749","2021-10-05 14:34:53.043751",
  "2021-10-05 14:34:53.043754","2021-10-05 14:34:53.043756","2021-10-05 14:34:53.043759","2021-10-05 14:34:53.043761",
  "2021-10-05 14:34:53.043764","2021-10-05 14:34:53.043766","2021-10-05 14:34:53.043768","2021-10-05 14:34:53.043771",
  "2021-10-05 14:34:53.043774","2021-10-05 14:34:53.043776","2021-10-05 14:34:53.043779","2021-10-05 14:34:53.043781",
  "2021-10-05 14:34:53.043783","2021-10-05 14:34:53.043786","2021-10-05 14:34:53.043789","2021-10-05 14:34:53.043791",
  "2021-10-05 14:34:53.043793","2021-10-05 14:34:53.043796","2021-10-05 14:34:53.043798","2021-10-05 14:34:53.043801",
  "2021-10-05 14:34:53.043803","2021-10-05 14:34:53.043806","2021-10-05 14:34:53.043809","2021-10-05 14:34:53.043811",
  "2021-10-05 14:34:53.043813","2021-10-05 14:34:53.043816","2021-10-05 14:34:53.043818","2021-10-05 14:34:53.043821",
  "2021-10-05 14:34:53.043823","2021-10-05 14:34:53.043826","2021-10-05 14:34:53.043828","2021-10-05 14:34:53.043830",
  "2021-10-05 14:34:53.043834","2021-10-05 14:34:53.043836","2021-10-05 14:34:53.043839","2021-10-05 14:34:53.043841",
  "2021-10-05 14:34:53.043844","2021-10-05 14:34:53.043846","2021-10-05 14:34:53.043848","2021-10-05 14:34:53.043851",
  "2021-10-05 14:34:53.043853","2021-10-05 14:34:53.043856","2021-10-05 14:34:53.043858","2021-10-05 14:34:53.043861",
  "2021-10-05 14:34:53.043865","2021-10-05 14:34:53.043868","2021-10-05 14:34:53.043870","2021-10-05 14:34:53.043873",
  "2021-10-05 14:34:53.043875","2021-10-05 14:34:53.043878","2021-10-05 14:34:53.043880","2021-10-05 14:34:53.043883",
  "2021-10-05 14:34:53.043885","2021-10-05 14:34:53.043888","2021-10-05 14:34:53.043890","2021-10-05 14:34:53.043893",
  "2021-10-05 14:34:53.043895","2021-10-05 14:34:53.043897","2021-10-05 14:34:53.043900","2021-10-05 14:34:53.043902",
  "2021-10-05 14:34:53.043906","2021-10-05 14:34:53.043909","2021-10-05 14:34:53.043911","2021-10-05 14:34:53.043914",
  "2021-10-05 14:34:53.043916","2021-10-05 14:34:53.043919","2021-10-05 14:34:53.043922","2021-10-05 14:34:53.043924",
  "2021-10-05 14:34:53.043927","2021-10-05 14:34:53.043929","2021-10-05 14:34:53.043931","2021-10-05 14:34:53.043934",
  "2021-10-05 14:34:53.043936","2021-10-05 14:34:53.043939","2021-10-05 14:34:53.043941","2021-10-05 14:34:53.043943",
  "2021-10-05 14:34:53.043947","2021-10-05 14:34:53.043950","2021-10-05 14:34:53.043953","2021-10-05 14:34:53.043955",
  "2021-10-05 14:34:53.043957","2021-10-05 14:34:53.043960","2021-10-05 14:34:53.043962","2021-10-05 14:34:53.043964",
  "2021-10-05 14:34:53.043967","2021-10-05 14:34:53.043969","2021-10-05 14:34:53.043972","2021-10-05 14:34:53.043974",
  "2021-10-05 14:34:53.043976","2021-10-05 14:34:53.043979","2021-10-05 14:34:53.043981","2021-10-05 14:34:53.043984",
  "2021-10-05 14:34:53.043986","2021-10-05 14:34:53.043988","2021-10-05 14:34:53.043991","2021-10-05 14:34:53.043993",
  "2021-10-05 14:34:53.043997","2021-10-05 14:34:53.044000","2021-10-05 14:34:53.044002","2021-10-05 14:34:53.044005",
  "2021-10-05 14:34:53.044007","2021-10-05 14:34:53.044009","2021-10-05 14:34:53.044012","2021-10-05 14:34:53.044014",
  "2021-10-05 14:34:53.044017","2021-10-05 14:34:53.044019","2021-10-05 14:34:53.044021","2021-10-05 14:34:53.044024",
  "2021-10-05 14:34:53.044028","2021-10-05 14:34:53.044032","2021-10-05 14:34:53.044036","2021-10-05 14:34:53.044040",
  "2021-10-05 14:34:53.044044","2021-10-05 14:34:53.044048","2021-10-05 14:34:53.044050","2021-10-05 14:34:53.044053",
  "2021-10-05 14:34:53.044056","2021-10-05 14:34:53.044061","2021-10-05 14:34:53.044064"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
  $row = $i + 2
  $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}
